# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E, rows 16-18) on the "Hoja1" account-statement
# sheet lists the arrears period for each worker. All three rows currently
# show period "2507" (Jul-2025); roll it forward to "2508" (Aug-2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "2507"
$newValue = "2508"

$periodoRange = $ws.Range("E16:E18")
foreach ($cell in $periodoRange.Cells) {
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
